$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: touch a cell so it is materialised in the sheet XML as an empty
# <c r=".."/> element without altering its style (setting the font name to
# the workbook's existing default font is a style-preserving no-op).
function New-EmptyCell($addr) {
    $ws.Range($addr).Font.Name = "等线"
}

# --- Existing cells whose string values change --------------------------
$ws.Range("I1").Value = "i1"
$ws.Range("I2").Value = "i2"
$ws.Range("H3").Value = "x var"
$ws.Range("I3").Value = "{i3}"
$ws.Range("C4").Value = "[A1]"
$ws.Range("I5").Value = "move to i5"

# --- Row 1: add J1 (empty) ----------------------------------------------
New-EmptyCell "J1"

# --- Row 2: add J2 (empty) ----------------------------------------------
New-EmptyCell "J2"

# --- Row 3: add J3 (empty) ----------------------------------------------
New-EmptyCell "J3"

# --- Row 4: add J4 (empty); F4:I4 already exist --------------------------
New-EmptyCell "J4"

# --- Row 5: add J5 (empty); F5:H5 already exist ---------------------------
New-EmptyCell "J5"

# --- Rows 6-8: fully empty rows across A,B,C,F,G,H,I,J --------------------
foreach ($r in 6..8) {
    foreach ($col in @("A","B","C","F","G","H","I","J")) {
        New-EmptyCell "$col$r"
    }
}

# --- Row 9 ----------------------------------------------------------------
$ws.Range("A9").Value = "a9"
foreach ($col in @("B","C","F","G","H","I","J")) {
    New-EmptyCell "$col`9"
}

# --- Row 10 -----------------------------------------------------------
$ws.Range("A10").Value = "a10"
$ws.Range("B10").Value = "x"
foreach ($col in @("C","F","G","H","I","J")) {
    New-EmptyCell "$col`10"
}

# --- Row 11 (sparse: only A, B, I, J) --------------------------------------
New-EmptyCell "A11"
$ws.Range("B11").Value = "xx"
New-EmptyCell "I11"
New-EmptyCell "J11"

# --- Row 12 -----------------------------------------------------------
$ws.Range("A12").Value = "{a11}"
foreach ($col in @("B","C","F","G","H","I","J")) {
    New-EmptyCell "$col`12"
}

# --- Selection moves to F2 -------------------------------------------------
$ws.Range("F2").Select() | Out-Null
